$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.848.89'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '2.910.03'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D6").Value = '''146.52'
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '2.907.78'
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").Value = '''6.84'
$ws.Range("E10").Value = '  -7.16%  '
$ws.Range("D11").Value = '''0.150'
$ws.Range("E11").Value = '  +4.97%  '
$ws.Range("E12").Value = '  -2.81%  '
$ws.Range("D13").Value = '''0.0000237'
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").Value = '''32.81'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '3.392.36'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("D17").Value = '61.866.21'
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Value = '2.911.33'
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("D20").Value = '''436.71'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").Value = '''13.40'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '''0.660'
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("D23").Value = '''6.96'
$ws.Range("E23").Value = '  -1.98%  '
$ws.Range("D24").Value = '''81.05'
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").Value = '''11.97'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = '''10.25'
$ws.Range("E26").Value = '  -8.09%  '
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = '''0.0000107'
$ws.Range("E29").Value = '  +22.58%  '
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("E33").Value = '  +0.37%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("D36").Value = '''0.973'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("E38").Value = '  +4.44%  '
$ws.Range("D39").Value = '''49.18'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("E40").Value = '  -2.63%  '
$ws.Range("D41").Value = '''8.38'
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").Value = '''0.273'
$ws.Range("E43").Value = '  -3.19%  '
$ws.Range("D44").Value = '''39.00'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = '2.701.34'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = '''133.98'
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").Value = '''342.61'
$ws.Range("E48").Value = '  -5.30%  '
$ws.Range("D50").Value = '''0.103'
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("D51").Value = '''22.31'
$ws.Range("E51").Value = '  -2.58%  '
